$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item(1,1).End(-4121).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text
    if ($text -ne $null -and $text -ne "") {
        $parts = $text.Split(",")
        if ($parts.Length -ge 2) {
            $first = $parts[0].Trim()
            $second = $parts[1].Trim()
            if ($second -eq "System" -and ($first -eq "dnasr281@gmail.com" -or $first -eq "backup@backdoor.com")) {
                $rest = ""
                for ($i = 2; $i -lt $parts.Length; $i++) {
                    $rest = $rest + ", " + $parts[$i].Trim()
                }
                $newText = "System, " + $first + $rest
                $cell.Value = $newText
            }
        }
    }
}
